$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their text representation
# (values such as "69.235.78" or "1.00" must stay literal strings,
# matching the workbook author's original inline-string cells,
# rather than being auto-coerced into numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 2).Value = "Bitcoin"
$ws.Cells.Item(2, 4).Value = "69.235.78"
$ws.Cells.Item(2, 5).Value = "  +5.68%  "

$ws.Cells.Item(3, 2).Value = "Ethereum"
$ws.Cells.Item(3, 4).Value = "3.558.20"
$ws.Cells.Item(3, 5).Value = "  +4.90%  "

$ws.Cells.Item(4, 2).Value = "TetherUSD"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.11%  "

$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 4).Value = "589.16"
$ws.Cells.Item(5, 5).Value = "  +5.18%  "

$ws.Cells.Item(6, 2).Value = "Solana"
$ws.Cells.Item(6, 4).Value = "191.77"
$ws.Cells.Item(6, 5).Value = "  +8.86%  "

$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 4).Value = "0.642"
$ws.Cells.Item(7, 5).Value = "  +1.74%  "

$ws.Cells.Item(8, 2).Value = "LidoStakedEther"
$ws.Cells.Item(8, 4).Value = "3.550.46"
$ws.Cells.Item(8, 5).Value = "  +4.91%  "

$ws.Cells.Item(9, 2).Value = "USDC"
$ws.Cells.Item(9, 4).Value = "1.00"
$ws.Cells.Item(9, 5).Value = "  -0.06%  "

$ws.Cells.Item(10, 2).Value = "Dogecoin"
$ws.Cells.Item(10, 4).Value = "0.182"
$ws.Cells.Item(10, 5).Value = "  +4.64%  "

$ws.Cells.Item(11, 2).Value = "Cardano"
$ws.Cells.Item(11, 4).Value = "0.659"
$ws.Cells.Item(11, 5).Value = "  +3.61%  "

$ws.Cells.Item(12, 2).Value = "Avalanche"
$ws.Cells.Item(12, 4).Value = "58.13"
$ws.Cells.Item(12, 5).Value = "  +8.70%  "

$ws.Cells.Item(13, 2).Value = "ShibaInu"
$ws.Cells.Item(13, 4).Value = "0.0000293"
$ws.Cells.Item(13, 5).Value = "  +5.49%  "

$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 4).Value = "9.64"
$ws.Cells.Item(14, 5).Value = "  +4.77%  "

$ws.Cells.Item(15, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(15, 4).Value = "4.124.74"
$ws.Cells.Item(15, 5).Value = "  +4.76%  "

$ws.Cells.Item(16, 2).Value = "Chainlink"
$ws.Cells.Item(16, 4).Value = "19.20"
$ws.Cells.Item(16, 5).Value = "  +5.01%  "

$ws.Cells.Item(17, 2).Value = "WrappedEther"
$ws.Cells.Item(17, 4).Value = "3.563.06"
$ws.Cells.Item(17, 5).Value = "  +4.89%  "

$ws.Cells.Item(18, 2).Value = "WrappedBTC"
$ws.Cells.Item(18, 4).Value = "69.281.22"
$ws.Cells.Item(18, 5).Value = "  +5.65%  "

$ws.Cells.Item(19, 2).Value = "Uniswap"
$ws.Cells.Item(19, 4).Value = "12.40"
$ws.Cells.Item(19, 5).Value = "  +4.76%  "

$ws.Cells.Item(20, 2).Value = "TRON"
$ws.Cells.Item(20, 4).Value = "0.120"
$ws.Cells.Item(20, 5).Value = "  +0.59%  "

$ws.Cells.Item(21, 2).Value = "Polygon"
$ws.Cells.Item(21, 4).Value = "1.04"
$ws.Cells.Item(21, 5).Value = "  +4.08%  "

$ws.Cells.Item(22, 2).Value = "BitcoinCash"
$ws.Cells.Item(22, 4).Value = "504.43"
$ws.Cells.Item(22, 5).Value = "  +4.77%  "

$ws.Cells.Item(23, 2).Value = "Toncoin"
$ws.Cells.Item(23, 4).Value = "5.51"
$ws.Cells.Item(23, 5).Value = "  +11.77%  "

$ws.Cells.Item(24, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(24, 4).Value = "17.03"
$ws.Cells.Item(24, 5).Value = "  +18.91%  "

$ws.Cells.Item(25, 2).Value = "PancakeSwap"
$ws.Cells.Item(25, 4).Value = "4.44"
$ws.Cells.Item(25, 5).Value = "  +8.17%  "

$ws.Cells.Item(26, 2).Value = "Litecoin"
$ws.Cells.Item(26, 4).Value = "91.01"
$ws.Cells.Item(26, 5).Value = "  +1.45%  "

$ws.Cells.Item(27, 2).Value = "ImmutableX"
$ws.Cells.Item(27, 4).Value = "3.03"
$ws.Cells.Item(27, 5).Value = "  +3.68%  "

$ws.Cells.Item(28, 2).Value = "RenderToken"
$ws.Cells.Item(28, 4).Value = "11.16"
$ws.Cells.Item(28, 5).Value = "  +4.80%  "

$ws.Cells.Item(29, 2).Value = "Filecoin"
$ws.Cells.Item(29, 4).Value = "9.25"
$ws.Cells.Item(29, 5).Value = "  +6.04%  "

$ws.Cells.Item(30, 2).Value = "EthereumClassic"
$ws.Cells.Item(30, 4).Value = "31.89"
$ws.Cells.Item(30, 5).Value = "  +1.89%  "

$ws.Cells.Item(31, 2).Value = "NEARProtocol"
$ws.Cells.Item(31, 4).Value = "7.46"
$ws.Cells.Item(31, 5).Value = "  +14.01%  "

$ws.Cells.Item(32, 2).Value = "Bittensor"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(32, 4).Value = "616.32"
$ws.Cells.Item(32, 5).Value = "  +7.19%  "

$ws.Cells.Item(33, 2).Value = "Cosmos"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(33, 4).Value = "12.11"
$ws.Cells.Item(33, 5).Value = "  +5.27%  "

$ws.Cells.Item(34, 2).Value = "OKB"
$ws.Cells.Item(34, 4).Value = "65.25"
$ws.Cells.Item(34, 5).Value = "  +3.93%  "

$ws.Cells.Item(35, 2).Value = "Hedera"
$ws.Cells.Item(35, 4).Value = "0.114"
$ws.Cells.Item(35, 5).Value = "  +5.98%  "

$ws.Cells.Item(36, 2).Value = "PEPE"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(36, 4).Value = "0.0₃0822"
$ws.Cells.Item(36, 5).Value = "  +11.35%  "

$ws.Cells.Item(37, 2).Value = "Kaspa"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(37, 4).Value = "0.148"
$ws.Cells.Item(37, 5).Value = "  +4.94%  "

$ws.Cells.Item(38, 2).Value = "Dai"
$ws.Cells.Item(38, 4).Value = "1.00"
$ws.Cells.Item(38, 5).Value = "  -0.02%  "

$ws.Cells.Item(39, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(39, 4).Value = "37.74"
$ws.Cells.Item(39, 5).Value = "  +5.19%  "

$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(40, 4).Value = "3.64"
$ws.Cells.Item(40, 5).Value = "  +0.69%  "

$ws.Cells.Item(41, 2).Value = "TheGraph"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(41, 4).Value = "0.396"
$ws.Cells.Item(41, 5).Value = "  +5.68%  "

$ws.Cells.Item(42, 2).Value = "Maker"
$ws.Cells.Item(42, 4).Value = "3.314.39"
$ws.Cells.Item(42, 5).Value = "  +6.95%  "

$ws.Cells.Item(43, 2).Value = "ThetaToken"
$ws.Cells.Item(43, 4).Value = "3.05"
$ws.Cells.Item(43, 5).Value = "  +9.02%  "

$ws.Cells.Item(44, 2).Value = "Fetch.AI"
$ws.Cells.Item(44, 4).Value = "2.69"
$ws.Cells.Item(44, 5).Value = "  +10.62%  "

$ws.Cells.Item(45, 2).Value = "VeChain"
$ws.Cells.Item(45, 4).Value = "0.0441"
$ws.Cells.Item(45, 5).Value = "  +5.57%  "

$ws.Cells.Item(46, 2).Value = "dogwifhat"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(46, 4).Value = "2.90"
$ws.Cells.Item(46, 5).Value = "  +22.70%  "

$ws.Cells.Item(47, 2).Value = "ApeXProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(47, 4).Value = "3.27"
$ws.Cells.Item(47, 5).Value = "  +3.22%  "

$ws.Cells.Item(48, 2).Value = "Stellar"
$ws.Cells.Item(48, 4).Value = "0.137"
$ws.Cells.Item(48, 5).Value = "  +2.07%  "

$ws.Cells.Item(49, 2).Value = "THORChain"
$ws.Cells.Item(49, 4).Value = "9.06"
$ws.Cells.Item(49, 5).Value = "  +7.30%  "

$ws.Cells.Item(50, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(50, 4).Value = "1.00"
$ws.Cells.Item(50, 5).Value = "  +0.07%  "

$ws.Cells.Item(51, 2).Value = "LidoDAOToken"
$ws.Cells.Item(51, 4).Value = "3.23"
$ws.Cells.Item(51, 5).Value = "  +3.72%  "
